$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, shifting existing rows 11..93 down to 12..94.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with the new data record.
$ws.Range("A11").Value = 4
$ws.Range("B11").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C11").Value = "Los Lagos"
$ws.Range("D11").Value = 44532
$ws.Range("E11").Value = 10
$ws.Range("F11").Value = 100112022
$ws.Range("G11").Value = "Arveja Verde"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 35
$ws.Range("K11").Value = 18000
$ws.Range("L11").Value = 18000
$ws.Range("M11").Value = 18000
$ws.Range("N11").Value = "$/saco 25 kilos"
$ws.Range("O11").Value = "Región del Maule"
$ws.Range("P11").Value = 720
$ws.Range("Q11").Value = 25
$ws.Range("R11").Value = "Hortaliza"

# Apply the date number format (style index 2 in the original file) to D11,
# matching the other date cells in column D.
$ws.Range("D11").NumberFormat = $ws.Range("D12").NumberFormat
